$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "etc" example table (rows 59-62): same First/Last/age/gender sample
# data used earlier in the sheet, plus a new "etc" column of numbers.
# ---------------------------------------------------------------------------
$ws.Range("B59").Value = "First Name"
$ws.Range("C59").Value = "Last name"
$ws.Range("D59").Value = "age"
$ws.Range("E59").Value = "gender"
$ws.Range("F59").Value = "etc"

$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Roddy"
$ws.Range("C60").Value = "Wiliams"
$ws.Range("D60").Value = 34
$ws.Range("E60").Value = "male"
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 3

$ws.Range("A61").Value = 2
$ws.Range("B61").Value = "Max"
$ws.Range("C61").Value = "Tiff"
$ws.Range("D61").Value = 74
$ws.Range("E61").Value = "male"
$ws.Range("F61").Value = 4
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 6

$ws.Range("A62").Value = 3
$ws.Range("B62").Value = "Lili"
$ws.Range("C62").Value = "Abrams"
$ws.Range("D62").Value = 23
$ws.Range("E62").Value = "female"
$ws.Range("F62").Value = 10
$ws.Range("G62").Value = 11
$ws.Range("H62").Value = 12

# ---------------------------------------------------------------------------
# Transposed rendition of the same table a few rows below (rows 67-71): each
# person's record becomes a column (B, E, H), headed by its row index, with
# the field names running down column A.
# ---------------------------------------------------------------------------
$ws.Range("B67").Value = 1
$ws.Range("E67").Value = 2
$ws.Range("H67").Value = 3

$ws.Range("A68").Value = "First Name"
$ws.Range("B68").Value = "Roddy"
$ws.Range("E68").Value = "Max"
$ws.Range("H68").Value = "Lili"

$ws.Range("A69").Value = "Last name"
$ws.Range("B69").Value = "Wiliams"
$ws.Range("E69").Value = "Tiff"
$ws.Range("H69").Value = "Abrams"

$ws.Range("A70").Value = "age"
$ws.Range("B70").Value = 34
$ws.Range("E70").Value = 74
$ws.Range("H70").Value = 23

$ws.Range("A71").Value = "gender"
$ws.Range("B71").Value = "male"
$ws.Range("E71").Value = "male"
$ws.Range("H71").Value = "female"

# ---------------------------------------------------------------------------
# Reflect the scrolled-down view/selection state.
# ---------------------------------------------------------------------------
$ws.Range("B60").Select()
